$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "92.368.01"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.81%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.281.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.96%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "630.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.92%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.412"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +12.42%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.720"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.67%  "

# Row 9
$ws.Range("E9").Value = "  -0.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.276.89"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.98%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.588"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.54%  "

# Row 12
$ws.Range("E12").Value = "  +0.03%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000263"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.18%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.44%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.883.97"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.71%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "92.171.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.31%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.310.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.65%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.56%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "439.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.77%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.88"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.55%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000192"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +49.80%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.28"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.20%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.39"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.18%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.27"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.31%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.472.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.54%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "76.95"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.00%  "

# Row 30
$ws.Range("E30").Value = "  +2.31%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.24%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.86%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "555.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.69%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.16"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.35%  "

# Row 35
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.93"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +24.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.58%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "22.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.60%  "

# Row 39
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.90%  "

# Row 40
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.45%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.393"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.67%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.66%  "

# Row 44
$ws.Range("E44").Value = "  +0.09%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "150.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.70%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "180.43"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.02%  "

# Row 47
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.46%  "

# Row 48
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.131"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.69%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.28"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.633"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.26%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
